$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven cell updates derived from the target diff.
# Each entry only lists the columns that actually change for that row.
$updates = @(
    @{ Row=2; D="328.28"; E="5.99%" },
    @{ Row=3; D="40.21"; E="8.10%" },
    @{ Row=4; D="5.659"; E="10.33%" },
    @{ Row=5; D="0.08132"; E="3.65%" },
    @{ Row=6; D="4.556"; E="3.41%" },
    @{ Row=7; D="8.694"; E="4.98%" },
    @{ Row=8; D="1.974"; E="5.27%" },
    @{ Row=9; E="0.42%" },
    @{ Row=10; D="0.9502"; E="2.65%" },
    @{ Row=11; D="0.1264"; E="7.07%" },
    @{ Row=12; D="0.1991"; E="5.02%" },
    @{ Row=13; D="0.09220"; E="3.52%" },
    @{ Row=14; D="0.03567" },
    @{ Row=15; D="0.09620"; E="0.27%" },
    @{ Row=16; D="0.001319"; E="-4.35%" },
    @{ Row=17; B="TigerCash"; C="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D="0.006017"; E="-2.99%" },
    @{ Row=18; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.376"; E="-0.40%" },
    @{ Row=19; B="BitpandaEcosystemToken"; C="https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D="0.3515"; E="1.65%" },
    @{ Row=20; B="MCDex"; C="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; D="7.504"; E="17.06%" },
    @{ Row=21; B="ProBitToken"; C="https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; D="0.1400"; E="8.35%" },
    @{ Row=22; B="ZBToken"; C="https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; D="0.2549"; E="5.95%" },
    @{ Row=23; B="CoinExToken"; C="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D="0.04444"; E="2.20%" },
    @{ Row=24; D="0.001252"; E="4.36%" },
    @{ Row=25; D="0.004320"; E="1.10%" },
    @{ Row=26; D="0.0001191"; E="-14.87%" },
    @{ Row=27; D="0.0003991"; E="37.48%" },
    @{ Row=39; D="0.02522"; E="17.07%" },
    @{ Row=40; D="0.05218"; E="4.43%" },
    @{ Row=41; D="0.007832"; E="3.50%" },
    @{ Row=42; D="0.1435"; E="6.16%" },
    @{ Row=43; D="0.009054"; E="6.62%" },
    @{ Row=44; D="0.002191"; E="8.91%" },
    @{ Row=45; D="0.01055"; E="33.81%" },
    @{ Row=46; D="0.00006714"; E="2.08%" },
    @{ Row=47; D="0.00000000750"; E="-0.08%" },
    @{ Row=48; D="0.002873"; E="-12.69%" },
    @{ Row=49; E="59.19%" },
    @{ Row=50; D="0.00002101"; E="-0.08%" },
    @{ Row=51; D="0.0002001"; E="-0.08%" },
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
